# "First start of renaming out -> outstanding"
# Rename the out_now / out_future_min / out_future_max headers to
# outstanding_now / outstanding_future_min / outstanding_future_max.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "outstanding_now"
$ws.Range("C1").Value = "outstanding_future_min"
$ws.Range("D1").Value = "outstanding_future_max"

# The longer header text no longer fits the old column widths, so widen
# (best-fit) columns B:D to match the new header labels.
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666
$ws.Columns.Item(3).ColumnWidth = 22.166666666666668
$ws.Columns.Item(4).ColumnWidth = 22.5

# Move the active selection, matching the author's cursor position after
# the edit.
$ws.Range("J14").Select() | Out-Null

$wb.Save()
